# "Exercises : day 1 complete"
# Tidy up the column-A variable names in the malaria dataset key sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the "enhanced_vegitation_index" typo -> "enhanced_vegetation_index"
$ws.Range("A32").Value = "enhanced_vegetation_index"

# "itnuse" is the type-0 net-usage variable; rename it to line up with the
# itnuse1 / itnuse2 / itnuse3 siblings already in the list.
$ws.Range("A22").Value = "itnuse0"

# Clarify that this column holds an index, not a raw suitability value.
$ws.Range("A30").Value = "mosquito_temperature_suitability_index"

# Reflect where the cursor ended up and resize column A for the longer names.
$null = $ws.Range("B9").Select()
$ws.Columns.Item(1).AutoFit() | Out-Null
